# Add input-cable activity lower bound block (rows 15-18), mirroring the
# existing ~TFM_INS block in rows 9-12 but for ACT_BND instead of CAP_BND.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header block (rows 9:10) formatting down to the new block (15:16)
$ws.Range("C9:Q10").Copy() | Out-Null
$ws.Range("C15:Q16").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 15: ~TFM_INS marker (copies the value along with the format)
$ws.Range("C15").Value = "~TFM_INS"

# Row 16: column headers
$ws.Range("C16").Value = "TimeSlice"
$ws.Range("D16").Value = "Year"
$ws.Range("E16").Value = "LimType"
$ws.Range("F16").Value = "Attribute"
$ws.Range("G16").Value = "DKE"
$ws.Range("H16").Value = "DKW"
$ws.Range("I16").Value = "DKISLBH"
$ws.Range("J16").Value = "DKISL1"
$ws.Range("K16").Value = "DKISL2"
$ws.Range("L16").Value = "DKISL3"
$ws.Range("M16").Value = "Pset_PN"
$ws.Range("N16").Value = "Cset_CN"

# Row 17: lower bound on input cable DKE-DKISLBH
$ws.Range("D17").Value = 2030
$ws.Range("E17").Value = "LO"
$ws.Range("M17").Value = "TB_ELCC_DKE_DKISLBH_01"

# Row 18: lower bound on input cable DKW-DKISL1
$ws.Range("D18").Value = 2030
$ws.Range("E18").Value = "LO"
$ws.Range("M18").Value = "TB_ELCC_DKW_DKISL1_01"

# Attribute + bound values added last so the shared-string table order
# matches (TB_ELCC_DKE_DKISLBH_01, TB_ELCC_DKW_DKISL1_01, ACT_BND)
$ws.Range("F17").Value = "ACT_BND"
$ws.Range("I17").Value = 10

$ws.Range("F18").Value = "ACT_BND"
$ws.Range("J18").Value = 10

# Restore the selection to where the user ended up after entering the data
$ws.Range("P28").Select() | Out-Null
